$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 253; everything currently at/after row 253
# (rows 253:337) shifts down by one (to 254:338), preserving all of its data.
$ws.Rows("253:253").Insert()

# Populate the newly inserted row 253 with the new weekly price record.
$ws.Range("A253").Value = 4
$ws.Range("B253").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C253").Value = "Los Lagos"
$ws.Range("D253").Value = 44876
$ws.Range("E253").Value = 10
$ws.Range("F253").Value = 100112021
$ws.Range("G253").Value = "Ají"
$ws.Range("H253").Value = "Inferno"
$ws.Range("I253").Value = "Segunda"
$ws.Range("J253").Value = 180
$ws.Range("K253").Value = 20000
$ws.Range("L253").Value = 20000
$ws.Range("M253").Value = 20000
$ws.Range("N253").Value = "$/caja 10 kilos"
$ws.Range("O253").Value = "Región de Arica y Parinacota"
$ws.Range("P253").Value = 2000
$ws.Range("Q253").Value = 10
$ws.Range("R253").Value = "Hortaliza"
